$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition) - rows 2-22 in column F get updated counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 113
$ws1.Range("F3").Value = 205
$ws1.Range("F5").Value = 6510
$ws1.Range("F7").Value = 426
$ws1.Range("F9").Value = 5871
$ws1.Range("F14").Value = 81
$ws1.Range("F15").Value = 380
$ws1.Range("F16").Value = 86
$ws1.Range("F18").Value = 332
$ws1.Range("F19").Value = 36
$ws1.Range("F21").Value = 4189
$ws1.Range("F22").Value = 32

# Sheet "全部类型" (all types) - same updates, but row 22 event is at row 23 here
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 113
$ws4.Range("F3").Value = 205
$ws4.Range("F5").Value = 6510
$ws4.Range("F7").Value = 426
$ws4.Range("F9").Value = 5871
$ws4.Range("F14").Value = 81
$ws4.Range("F15").Value = 380
$ws4.Range("F16").Value = 86
$ws4.Range("F18").Value = 332
$ws4.Range("F19").Value = 36
$ws4.Range("F21").Value = 4189
$ws4.Range("F23").Value = 32
